# Remove the "Appendix: Quick prototype" section in its entirety:
#   - the "Appendix: Quick prototype" Heading2 paragraph
#   - the blank spacer paragraph under it
#   - the three "Figure: PDF page N" reference paragraphs
#   - the three paragraphs holding the embedded prototype screenshots
#   - the blank spacer paragraph that follows the last image
# The surviving "Appendix: Links" Heading2 section (and everything above
# the prototype appendix) is left untouched.

$d = $word.ActiveDocument

# Locate the start of the block to remove: the "Appendix: Quick prototype"
# heading paragraph.
$startRange = $d.Content
$found = $startRange.Find.Execute("Appendix: Quick prototype", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'Appendix: Quick prototype' heading"
}
$deleteStart = $startRange.Paragraphs.Item(1).Range.Start

# Locate the end of the block to remove: the second occurrence of the text
# "Appendix: Links" in the document is the Heading2 that must be kept, so
# the block we delete stops right before it. (The first occurrence is just
# plain text inside the earlier "Appendix: Links" paragraph higher up.)
$endRange = $d.Content
$null = $endRange.Find.Execute("Appendix: Links", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$found2 = $endRange.Find.Execute("Appendix: Links", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find the surviving 'Appendix: Links' heading"
}
$deleteEnd = $endRange.Paragraphs.Item(1).Range.Start

# Delete the whole prototype-appendix block (heading, spacer paragraphs,
# figure captions, and embedded images) in one shot.
$d.Range($deleteStart, $deleteEnd).Delete()
